$d = $word.ActiveDocument

# --- Edit 1: paragraph "(iv) Acceptable turnaround time(s) ..." (List3 style) ---
# Split the single run into two runs: "(iv)" and " Acceptable turnaround ... met."
# This paragraph consists of exactly one run, so we rebuild the whole paragraph
# (preserving its original pPr / rsid attributes) to keep the pStyle intact.
$oldText1 = "(iv) Acceptable turnaround time(s) for warranty corrective actions taken by the contractor should be specified, and consideration should be given to using liquidated damages or charging the contractor for product replacement costs when specified turnaround times are not met."

$rng1 = $d.Content
$found1 = $rng1.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Edit 1: could not find the '(iv) Acceptable turnaround...' paragraph text"
}
$rng1.Text = ""
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00537A74" w:rsidRDefault="00537A74" w:rsidP="00537A74"><w:pPr><w:pStyle w:val="List3"/></w:pPr><w:r><w:t>(iv)</w:t></w:r><w:r><w:t xml:space="preserve"> Acceptable turnaround time(s) for warranty corrective actions taken by the contractor should be specified, and consideration should be given to using liquidated damages or charging the contractor for product replacement costs when specified turnaround times are not met.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $rng1.InsertXML($xml1)

# --- Edit 2: paragraph "(5) Markings. ..." (List2 style) ---
# Split the leading run "(5) " into two runs: "(5)" and " ".
# This paragraph has several other runs after it (Markings, the body text, a
# noBreakHyphen run, ...) so we only touch the small sub-range holding "(5) ",
# leaving the rest of the paragraph (and its pPr) untouched.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("(5) Markings", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Edit 2: could not find the '(5) Markings' paragraph text"
}
$start2 = $rng2.Start
$target2 = $d.Range($start2, $start2 + 4)
if ($target2.Text -ne "(5) ") {
    throw "Edit 2: unexpected text at target range: [$($target2.Text)]"
}
$target2.Text = ""
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>(5)</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $target2.InsertXML($xml2)
